$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 72
$ws.Range("I9").Value = 72
$ws.Range("K9").Value = 72
$ws.Range("M9").Value = 97
$ws.Range("H64").Value = 4999.8887
$ws.Range("I64").Value = 4999
$ws.Range("K64").Value = 4999
$ws.Range("M64").Value = -4751
$ws.Range("H67").Value = 4999.8887
$ws.Range("I67").Value = 4999
$ws.Range("K67").Value = 4999
$ws.Range("M67").Value = -4141
$ws.Range("H92").Value = 1210.08
$ws.Range("I92").Value = 1135.4706
$ws.Range("J92").Value = 1368.625
$ws.Range("K92").Value = 1135.4706
$ws.Range("L92").Value = 1368.625
$ws.Range("M92").Value = 112.5293999999999
$ws.Range("N92").Value = -3864.625
$ws.Range("H98").Value = 62501496
$ws.Range("I98").Value = 62501496
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 62501496
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -62499998
$ws.Range("H122").Value = 62501496
$ws.Range("I122").Value = 62501496
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 187504488
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -187502038
$ws.Range("H137").Value = 4512.5
$ws.Range("I137").Value = 2939.5
$ws.Range("K137").Value = 8818.5
$ws.Range("M137").Value = -6268.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1955.5555
$ws.Range("I2").Value = 1764.5
$ws.Range("K2").Value = 1764.5
$ws.Range("M2").Value = -1651.5
$ws.Range("H32").Value = 20835462
$ws.Range("I32").Value = 22728776
$ws.Range("K32").Value = 22728776
$ws.Range("M32").Value = -22728489
$ws.Range("H61").Value = 35789604
$ws.Range("I61").Value = 71430056
$ws.Range("J61").Value = 149153.42
$ws.Range("K61").Value = 71430056
$ws.Range("L61").Value = 149153.42
$ws.Range("M61").Value = -71429844
$ws.Range("N61").Value = -149577.42
$ws.Range("H74").Value = 11374758
$ws.Range("I74").Value = 22729000
$ws.Range("J74").Value = 20515.182
$ws.Range("K74").Value = 22729000
$ws.Range("L74").Value = 20515.182
$ws.Range("M74").Value = -22728126
$ws.Range("N74").Value = -22263.182
$ws.Range("H77").Value = 11374758
$ws.Range("I77").Value = 22729000
$ws.Range("J77").Value = 20515.182
$ws.Range("K77").Value = 113645000
$ws.Range("L77").Value = 102575.91
$ws.Range("M77").Value = -113640632
$ws.Range("N77").Value = -111311.91
$ws.Range("H88").Value = 1520.4706
$ws.Range("I88").Value = 1375.1111
$ws.Range("K88").Value = 1375.1111
$ws.Range("M88").Value = -969.1111000000001
$ws.Range("H91").Value = 1520.4706
$ws.Range("I91").Value = 1375.1111
$ws.Range("K91").Value = 1375.1111
$ws.Range("M91").Value = 28.88889999999992
$ws.Range("H116").Value = 1955.5555
$ws.Range("I116").Value = 1764.5
$ws.Range("K116").Value = 1764.5
$ws.Range("M116").Value = 529.5
$ws.Range("H132").Value = 4957.8687
$ws.Range("I132").Value = 2683.3572
$ws.Range("K132").Value = 8050.071599999999
$ws.Range("M132").Value = -5520.071599999999
$ws.Range("H136").Value = 35789604
$ws.Range("I136").Value = 71430056
$ws.Range("J136").Value = 149153.42
$ws.Range("K136").Value = 214290168
$ws.Range("L136").Value = 447460.26
$ws.Range("M136").Value = -214287618
$ws.Range("N136").Value = -452560.26
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1955.5555
$ws.Range("I3").Value = 1764.5
$ws.Range("K3").Value = 1764.5
$ws.Range("M3").Value = -1650.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1168398.5
$ws.Range("I31").Value = 1780.125
$ws.Range("J31").Value = 1946144.1
$ws.Range("K31").Value = 1780.125
$ws.Range("L31").Value = 1946144.1
$ws.Range("M31").Value = -1485.125
$ws.Range("N31").Value = -1946734.1
$ws.Range("H34").Value = 1168398.5
$ws.Range("I34").Value = 1780.125
$ws.Range("J34").Value = 1946144.1
$ws.Range("K34").Value = 1780.125
$ws.Range("L34").Value = 1946144.1
$ws.Range("M34").Value = -1578.125
$ws.Range("N34").Value = -1946548.1
$ws.Range("H58").Value = 1843.3636
$ws.Range("I58").Value = 1712.8572
$ws.Range("K58").Value = 1712.8572
$ws.Range("M58").Value = -1509.8572
$ws.Range("H132").Value = 1821.4117
$ws.Range("J132").Value = 3740.6667
$ws.Range("L132").Value = 11222.0001
$ws.Range("N132").Value = -16282.0001
$ws.Range("H134").Value = 1434472.2
$ws.Range("I134").Value = 5000900
$ws.Range("K134").Value = 15002700
$ws.Range("M134").Value = -15000165
$ws.Range("H136").Value = 1843.3636
$ws.Range("I136").Value = 1712.8572
$ws.Range("K136").Value = 5138.571599999999
$ws.Range("M136").Value = -2588.571599999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 263.33334
$ws.Range("J17").Value = 263.33334
$ws.Range("L17").Value = 790.0000200000001
$ws.Range("N17").Value = -1128.00002
$ws.Range("H69").Value = 4435
$ws.Range("J69").Value = 4435
$ws.Range("L69").Value = 13305
$ws.Range("N69").Value = -14927
$ws.Range("H72").Value = 4435
$ws.Range("J72").Value = 4435
$ws.Range("L72").Value = 39915
$ws.Range("N72").Value = -48027
$ws.Range("H113").Value = 1836.8334
$ws.Range("I113").Value = 394.25
$ws.Range("K113").Value = 1182.75
$ws.Range("M113").Value = 987.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3684.6155
$ws.Range("I113").Value = 2983.3333
$ws.Range("K113").Value = 2983.3333
$ws.Range("M113").Value = -813.3332999999998
$ws.Range("H122").Value = 3503.7
$ws.Range("I122").Value = 3115.3333
$ws.Range("J122").Value = 6999
$ws.Range("K122").Value = 9345.999899999999
$ws.Range("L122").Value = 20997
$ws.Range("M122").Value = -6895.999899999999
$ws.Range("N122").Value = -25897
$ws.Range("H132").Value = 27780848
$ws.Range("I132").Value = 32261020
$ws.Range("J132").Value = 3780
$ws.Range("K132").Value = 96783060
$ws.Range("L132").Value = 11340
$ws.Range("M132").Value = -96780530
$ws.Range("N132").Value = -16400
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 51783.81
$ws.Range("I7").Value = 3527.5625
$ws.Range("K7").Value = 3527.5625
$ws.Range("M7").Value = -3415.5625
$ws.Range("H40").Value = 4363.2856
$ws.Range("I40").Value = 3454.111
$ws.Range("K40").Value = 3454.111
$ws.Range("M40").Value = -3318.111
$ws.Range("H61").Value = 9479.799999999999
$ws.Range("J61").Value = 39399
$ws.Range("L61").Value = 39399
$ws.Range("N61").Value = -39803
$ws.Range("H113").Value = 9479.799999999999
$ws.Range("J113").Value = 39399
$ws.Range("L113").Value = 39399
$ws.Range("N113").Value = -43739
$ws.Range("H122").Value = 5260.222
$ws.Range("I122").Value = 4633.483
$ws.Range("K122").Value = 13900.449
$ws.Range("M122").Value = -11450.449
$ws.Range("H126").Value = 51783.81
$ws.Range("I126").Value = 3527.5625
$ws.Range("K126").Value = 10582.6875
$ws.Range("M126").Value = -8112.6875
$ws.Range("H136").Value = 181373.67
$ws.Range("J136").Value = 165711.42
$ws.Range("L136").Value = 497134.26
$ws.Range("N136").Value = -502234.26
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 40000
$ws.Range("I7").Value = 40000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 40000
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -39887
$ws.Range("H9").Value = 12000
$ws.Range("J9").Value = 12000
$ws.Range("L9").Value = 12000
$ws.Range("N9").Value = -12280
$ws.Range("H64").Value = 64113.5
$ws.Range("J64").Value = 64113.5
$ws.Range("L64").Value = 64113.5
$ws.Range("N64").Value = -64609.5
$ws.Range("H67").Value = 64113.5
$ws.Range("J67").Value = 64113.5
$ws.Range("L67").Value = 64113.5
$ws.Range("N67").Value = -65829.5
$ws.Range("H107").Value = 23810556
$ws.Range("I107").Value = 33334482
$ws.Range("K107").Value = 100003446
$ws.Range("M107").Value = -100001526
$ws.Range("H110").Value = 34995
$ws.Range("J110").Value = 34995
$ws.Range("L110").Value = 34995
$ws.Range("N110").Value = -43175
$ws.Range("H132").Value = 4174.1177
$ws.Range("I132").Value = 4106.5386
$ws.Range("K132").Value = 12319.6158
$ws.Range("M132").Value = -9789.6158
$ws.Range("H136").Value = 1618.8422
$ws.Range("I136").Value = 1547.375
$ws.Range("K136").Value = 4642.125
$ws.Range("M136").Value = -2092.125

Write-Host "Applied 220 cell updates across 8 sheets."